# Implements: add a new "cardDetails" worksheet (sheetId 3) after the
# existing sheets, populate it with card-detail sample data, format the
# card-number column, and make it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet at the end of the workbook --------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "cardDetails"

# --- Header row ------------------------------------------------------------
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "card number"
$ws.Range("C1").Value = "CVS"
$ws.Range("D1").Value = "MM"
$ws.Range("E1").Value = "YY"

# --- Data row ----------------------------------------------------------
$ws.Range("A2").Value = "Ankit"
$ws.Range("B2").Value = 9652178523256980
$ws.Range("C2").Value = 652
$ws.Range("D2").Value = 11
$ws.Range("E2").Value = 2030

# --- Formatting: widen the "card number" column and give it a numeric
#     format (2 decimal places) ------------------------------------------
$cardNumberColumn = $ws.Columns.Item(2)
$cardNumberColumn.ColumnWidth = 29.36
$cardNumberColumn.NumberFormat = "0.00"

# --- Selection / active sheet -------------------------------------------
$ws.Range("B2").Select()
